# Update Ligand/Receptor/Edge metrics for F2-F2r sheet with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 9,16
$row = @(1.0, 0.3333333333333333, 0.04339299999999999, 0.130179, 0.0698021577815419, 0.0698021577815419, 3.0, 1.0, 4.536281, 13.608843, 0.0729209419042404, 0.07292094190424041, 0.196842841433, 1.771585572897, 0.005090039092378438, 0.005090039092378439)
for ($j = 0; $j -lt 16; $j++) { $data[0, $j] = $row[$j] }
$row = @(1.0, 0.3333333333333333, 0.04339299999999999, 0.130179, 0.0698021577815419, 0.0698021577815419, 3.0, 1.0, 44.50790266666667, 133.523708, 0.7154667412877611, 0.7154667412877612, 1.931331420414666, 17.381982783732, 0.04994112236281392, 0.04994112236281393)
for ($j = 0; $j -lt 16; $j++) { $data[1, $j] = $row[$j] }
$row = @(1.0, 0.3333333333333333, 0.04339299999999999, 0.130179, 0.0698021577815419, 0.0698021577815419, 3.0, 1.0, 13.16402266666667, 39.492068, 0.2116123168079984, 0.2116123168079984, 0.5712264355746667, 5.141037920172, 0.01477099632634954, 0.01477099632634954)
for ($j = 0; $j -lt 16; $j++) { $data[2, $j] = $row[$j] }
$row = @(2.0, 0.6666666666666666, 0.4205383333333333, 1.261615, 0.6764796878879081, 0.6764796878879081, 3.0, 1.0, 4.536281, 13.608843, 0.0729209419042404, 0.07292094190424041, 1.907680051271666, 17.169120461445, 0.04932953601987283, 0.04932953601987283)
for ($j = 0; $j -lt 16; $j++) { $data[3, $j] = $row[$j] }
$row = @(2.0, 0.6666666666666666, 0.4205383333333333, 1.261615, 0.6764796878879081, 0.6764796878879081, 3.0, 1.0, 44.50790266666667, 133.523708, 0.7154667412877611, 0.7154667412877612, 18.71727920760222, 168.45551286842, 0.4839987178405233, 0.4839987178405235)
for ($j = 0; $j -lt 16; $j++) { $data[4, $j] = $row[$j] }
$row = @(2.0, 0.6666666666666666, 0.4205383333333333, 1.261615, 0.6764796878879081, 0.6764796878879081, 3.0, 1.0, 13.16402266666667, 39.492068, 0.2116123168079984, 0.2116123168079984, 5.535976152202222, 49.82378536982, 0.1431514340275119, 0.1431514340275119)
for ($j = 0; $j -lt 16; $j++) { $data[5, $j] = $row[$j] }
$row = @(2.0, 0.6666666666666666, 0.1577256666666667, 0.473177, 0.2537181543305499, 0.2537181543305499, 3.0, 1.0, 4.536281, 13.608843, 0.0729209419042404, 0.07292094190424041, 0.7154879449123334, 6.439391504211001, 0.01850136679198913, 0.01850136679198914)
for ($j = 0; $j -lt 16; $j++) { $data[6, $j] = $row[$j] }
$row = @(2.0, 0.6666666666666666, 0.1577256666666667, 0.473177, 0.2537181543305499, 0.2537181543305499, 3.0, 1.0, 44.50790266666667, 133.523708, 0.7154667412877611, 0.7154667412877612, 7.020038620035112, 63.180347580316, 0.1815269010844238, 0.1815269010844238)
for ($j = 0; $j -lt 16; $j++) { $data[7, $j] = $row[$j] }
$row = @(2.0, 0.6666666666666666, 0.1577256666666667, 0.473177, 0.2537181543305499, 0.2537181543305499, 3.0, 1.0, 13.16402266666667, 39.492068, 0.2116123168079984, 0.2116123168079984, 2.076304251115111, 18.686738260036, 0.05368988645413697, 0.05368988645413698)
for ($j = 0; $j -lt 16; $j++) { $data[8, $j] = $row[$j] }

$ws.Range("E2:T10").Value = $data
